# Update crypto price/volume figures per the GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage each new value as a text-producing formula ( ="literal" ) so Excel
# keeps these as strings (matching the source data, e.g. "94.217.57" / "  +1.57%  ")
# instead of silently coercing number-looking text into a float.
$ws.Range("D2").Formula = '="94.217.57"'
$ws.Range("E2").Formula = '="  +1.57%  "'
$ws.Range("D3").Formula = '="3.083.91"'
$ws.Range("E3").Formula = '="  -1.18%  "'
$ws.Range("E4").Formula = '="  -0.10%  "'
$ws.Range("D5").Formula = '="234.39"'
$ws.Range("E5").Formula = '="  -3.09%  "'
$ws.Range("D6").Formula = '="609.03"'
$ws.Range("E6").Formula = '="  -1.13%  "'
$ws.Range("E7").Formula = '="  -0.05%  "'
$ws.Range("E8").Formula = '="  -5.34%  "'
$ws.Range("E9").Formula = '="  -0.08%  "'
$ws.Range("D10").Formula = '="0.822"'
$ws.Range("E10").Formula = '="  +12.31%  "'
$ws.Range("D11").Formula = '="3.077.93"'
$ws.Range("E11").Formula = '="  -1.25%  "'
$ws.Range("E12").Formula = '="  -3.17%  "'
$ws.Range("D13").Formula = '="93.941.86"'
$ws.Range("E13").Formula = '="  +1.61%  "'
$ws.Range("E14").Formula = '="  -5.52%  "'
$ws.Range("D15").Formula = '="33.96"'
$ws.Range("E15").Formula = '="  -1.43%  "'
$ws.Range("E16").Formula = '="  -4.23%  "'
$ws.Range("D17").Formula = '="3.647.78"'
$ws.Range("E17").Formula = '="  -1.44%  "'
$ws.Range("D18").Formula = '="3.070.69"'
$ws.Range("E18").Formula = '="  -1.58%  "'
$ws.Range("D19").Formula = '="3.63"'
$ws.Range("E19").Formula = '="  -3.25%  "'
$ws.Range("D20").Formula = '="14.53"'
$ws.Range("E20").Formula = '="  -1.28%  "'
$ws.Range("D21").Formula = '="5.74"'
$ws.Range("E21").Formula = '="  -1.09%  "'
$ws.Range("D22").Formula = '="441.57"'
$ws.Range("E22").Formula = '="  -1.47%  "'
$ws.Range("D23").Formula = '="8.81"'
$ws.Range("E23").Formula = '="  -6.77%  "'
$ws.Range("E24").Formula = '="  -6.90%  "'
$ws.Range("E25").Formula = '="  +5.03%  "'
$ws.Range("D26").Formula = '="5.54"'
$ws.Range("E26").Formula = '="  -4.60%  "'
$ws.Range("D27").Formula = '="84.99"'
$ws.Range("E27").Formula = '="  -2.23%  "'
$ws.Range("D28").Formula = '="11.91"'
$ws.Range("E28").Formula = '="  +1.22%  "'
$ws.Range("D29").Formula = '="3.252.01"'
$ws.Range("D31").Formula = '="0.245"'
$ws.Range("E31").Formula = '="  +5.67%  "'
$ws.Range("D32").Formula = '="0.178"'
$ws.Range("E32").Formula = '="  +5.24%  "'
$ws.Range("D33").Formula = '="0.122"'
$ws.Range("E33").Formula = '="  -10.44%  "'
$ws.Range("D34").Formula = '="9.14"'
$ws.Range("E34").Formula = '="  -1.43%  "'
$ws.Range("E35").Formula = '="  -0.70%  "'
$ws.Range("E36").Formula = '="  -3.20%  "'
$ws.Range("E37").Formula = '="  -4.27%  "'
$ws.Range("D38").Formula = '="25.60"'
$ws.Range("E38").Formula = '="  -2.41%  "'
$ws.Range("D39").Formula = '="1.87"'
$ws.Range("E39").Formula = '="  -2.01%  "'
$ws.Range("E40").Formula = '="  +1.19%  "'
$ws.Range("D41").Formula = '="23.98"'
$ws.Range("E41").Formula = '="  +3.90%  "'
$ws.Range("E42").Formula = '="  -11.81%  "'
$ws.Range("D43").Formula = '="467.58"'
$ws.Range("E43").Formula = '="  -3.53%  "'
$ws.Range("E44").Formula = '="  -3.96%  "'
$ws.Range("D46").Formula = '="3.12"'
$ws.Range("E46").Formula = '="  -11.10%  "'
$ws.Range("D47").Formula = '="159.77"'
$ws.Range("E47").Formula = '="  -1.38%  "'
$ws.Range("D48").Formula = '="1.85"'
$ws.Range("E48").Formula = '="  -3.96%  "'
$ws.Range("E49").Formula = '="  -2.68%  "'
$ws.Range("D50").Formula = '="43.73"'
$ws.Range("E50").Formula = '="  -0.89%  "'
$ws.Range("E51").Formula = '="  -0.07%  "'

# Flatten the helper formulas down to plain literal text values (no stray
# style/number-format side effects, no leftover formulas) via copy / paste-values.
$touched = $ws.Range("D2:E51")
$touched.Copy()
$touched.PasteSpecial(-4163)
$excel.CutCopyMode = 0
